$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.358.36"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "1.652.58"
$ws.Range("E3").Value = "  -2.86%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'311.77"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "'0.3912"
$ws.Range("E7").Value = "  -1.85%  "
$ws.Range("E8").Value = "  -3.83%  "
$ws.Range("D9").Value = "'1.004"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.378"
$ws.Range("E10").Value = "  -5.99%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'50.21"
$ws.Range("E11").Value = "  -6.31%  "
$ws.Range("D12").Value = "'0.08556"
$ws.Range("E12").Value = "  -2.91%  "
$ws.Range("D13").Value = "'24.99"
$ws.Range("E13").Value = "  -4.80%  "
$ws.Range("E14").Value = "  -4.21%  "
$ws.Range("E15").Value = "  -2.83%  "
$ws.Range("D16").Value = "'7.607"
$ws.Range("E16").Value = "  -4.53%  "
$ws.Range("D17").Value = "1.657.47"
$ws.Range("E17").Value = "  -6.15%  "
$ws.Range("D18").Value = "'93.18"
$ws.Range("E18").Value = "  -2.59%  "
$ws.Range("D19").Value = "'0.06949"
$ws.Range("E19").Value = "  -3.32%  "
$ws.Range("D20").Value = "'21.07"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").Value = "'6.999"
$ws.Range("E21").Value = "  -4.40%  "
$ws.Range("D22").Value = "'1.005"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "'13.80"
$ws.Range("E23").Value = "  -4.18%  "
$ws.Range("D24").Value = "24.368.54"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("D25").Value = "'2.341"
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("D26").Value = "'2.778"
$ws.Range("E26").Value = "  -3.52%  "
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("D28").Value = "'158.89"
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("D29").Value = "'5.751"
$ws.Range("E29").Value = "  -5.27%  "
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("D31").Value = "'8.169"
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("D32").Value = "'2.505"
$ws.Range("E32").Value = "  +10.63%  "
$ws.Range("D33").Value = "1.836.31"
$ws.Range("E33").Value = "  -9.34%  "
$ws.Range("D34").Value = "'0.03009"
$ws.Range("E34").Value = "  -5.47%  "
$ws.Range("D35").Value = "'0.08079"
$ws.Range("E35").Value = "  -5.90%  "
$ws.Range("D36").Value = "'0.9971"
$ws.Range("E36").Value = "  -3.29%  "
$ws.Range("D37").Value = "'6.846"
$ws.Range("E37").Value = "  -5.93%  "
$ws.Range("D38").Value = "'0.2761"
$ws.Range("E38").Value = "  -3.03%  "
$ws.Range("D39").Value = "'0.09447"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "'1.491"
$ws.Range("E40").Value = "  +0.93%  "
$ws.Range("D41").Value = "'10.17"
$ws.Range("E41").Value = "  -5.07%  "
$ws.Range("D42").Value = "'0.7775"
$ws.Range("E42").Value = "  -6.61%  "
$ws.Range("D43").Value = "'13.32"
$ws.Range("E43").Value = "  -6.12%  "
$ws.Range("D44").Value = "'16.37"
$ws.Range("E44").Value = "  -6.71%  "
$ws.Range("D45").Value = "'2.552"
$ws.Range("E45").Value = "  -5.67%  "
$ws.Range("D46").Value = "'0.7005"
$ws.Range("E46").Value = "  -5.57%  "
$ws.Range("D47").Value = "'4.140"
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").Value = "'0.08557"
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("D50").Value = "'1.300"
$ws.Range("E50").Value = "  -5.31%  "
$ws.Range("D51").Value = "'136.42"
$ws.Range("E51").Value = "  -2.21%  "
